$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: "presenting_worthy" (C) and "additional_notes" (D).
# add_totals_col (row 2) and add_totals_row (row 3) are flagged as not
# presentation-worthy and deprecated in favor of adorn_totals.
$ws.Range("C1").Value = "presenting_worthy"
$ws.Range("C2").Value = "no"
$ws.Range("C3").Value = "no"

$ws.Range("D1").Value = "additional_notes"
$ws.Range("D2").Value = "deprecated, use adorn_totals"
$ws.Range("D3").Value = "deprecated, use adorn_totals"

# Match the bold header style used for the existing headers in row 1.
$ws.Range("C1:D1").Font.Bold = $true

# Size the new columns to fit their contents.
$ws.Columns("C:D").AutoFit()

# Portrait orientation for printing.
$ws.PageSetup.Orientation = 1

# Move the active selection, as left by the editor.
$ws.Range("C6").Select() | Out-Null
